$wb = $excel.ActiveWorkbook

# --- Analysis Results sheet ---
$ws = $wb.Worksheets.Item("Analysis Results")

# Column H width: 20.9 -> 13.2 (closest achievable via pixel-quantized ColumnWidth is 13.166667)
$ws.Columns.Item(8).ColumnWidth = 12.3

$ws.Range("C2").Value = 87.16
$ws.Range("H2").Value = 90
$ws.Range("R2").Value = "Frequency in description matches metadata (monthly)"
$ws.Range("C3").Value = 76.68000000000001
$ws.Range("H3").Value = 90
$ws.Range("R3").Value = "Frequency in description matches metadata (monthly)"
$ws.Range("Q4").Value = "No"
$ws.Range("R4").Value = "Frequency in description (none) does not match metadata (monthly)"
$ws.Range("C5").Value = 46.40000000000001
$ws.Range("E5").Value = "WHAT, WHY, ESCALATION"
$ws.Range("H5").Value = 80
$ws.Range("Q5").Value = "No"
$ws.Range("R5").Value = "Frequency in description (none) does not match metadata (ongoing)"
$ws.Range("C6").Value = 43.2
$ws.Range("E6").Value = "WHAT, WHY, ESCALATION"
$ws.Range("H6").Value = 80
$ws.Range("Q6").Value = "No"
$ws.Range("R6").Value = "Frequency in description (none) does not match metadata (daily)"
$ws.Range("Q7").Value = "No"
$ws.Range("R7").Value = "Frequency in description (none) does not match metadata (monthly)"
$ws.Range("Q8").Value = "No"
$ws.Range("R8").Value = "Frequency in description (none) does not match metadata (as needed)"
$ws.Range("Q9").Value = "No"
$ws.Range("R9").Value = "Frequency in description (none) does not match metadata (quarterly)"
$ws.Range("Q10").Value = "No"
$ws.Range("R10").Value = "Frequency in description (none) does not match metadata (weekly)"
$ws.Range("Q11").Value = "No"
$ws.Range("R11").Value = "Frequency in description (none) does not match metadata (monthly)"
$ws.Range("C12").Value = 29.84544
$ws.Range("H12").Value = 90
$ws.Range("R12").Value = "Frequency in description matches metadata (monthly)"
$ws.Range("C13").Value = 46.68
$ws.Range("H13").Value = 90
$ws.Range("R13").Value = "Frequency in description matches metadata (weekly)"
$ws.Range("C14").Value = 45.40000000000001
$ws.Range("H14").Value = 90
$ws.Range("R14").Value = "Frequency in description matches metadata (quarterly)"
$ws.Range("R15").Value = "Frequency in description (none) does not match metadata (as needed)"
$ws.Range("C16").Value = 30.22
$ws.Range("H16").Value = 85
$ws.Range("Q16").Value = "No"
$ws.Range("R16").Value = "Frequency in description (none) does not match metadata (within 3 days)"
$ws.Range("Q17").Value = "No"
$ws.Range("R17").Value = "Frequency in description (none) does not match metadata (ongoing)"
$ws.Range("Q18").Value = "No"
$ws.Range("R18").Value = "Frequency in description (none) does not match metadata (regulatory)"
$ws.Range("Q19").Value = "No"
$ws.Range("R19").Value = "Frequency in description (none) does not match metadata (ongoing)"
$ws.Range("Q20").Value = "No"
$ws.Range("R20").Value = "Frequency in description (none) does not match metadata (prompt)"
$ws.Range("Q21").Value = "No"
$ws.Range("R21").Value = "Frequency in description (none) does not match metadata (internal)"
$ws.Range("C23").Value = 52.68
$ws.Range("E23").Value = "WHEN, WHY"
$ws.Range("H23").Value = 0
$ws.Range("C27").Value = 93.18000000000001
$ws.Range("H27").Value = 90
$ws.Range("R27").Value = "Frequency in description matches metadata (daily)"
$ws.Range("C28").Value = 91.8411136
$ws.Range("H28").Value = 90
$ws.Range("R28").Value = "Frequency in description matches metadata (monthly, quarterly)"
$ws.Range("C29").Value = 89.08
$ws.Range("H29").Value = 90
$ws.Range("M29").Value = "high"
$ws.Range("O29").Value = "Yes"
$ws.Range("R29").Value = "Frequency in description matches metadata (weekly, ad hoc)"
$ws.Range("C30").Value = 83.263488
$ws.Range("H30").Value = 90
$ws.Range("R30").Value = "Frequency in description matches metadata (daily)"
$ws.Range("C31").Value = 95.26173217391305
$ws.Range("H31").Value = 90
$ws.Range("R31").Value = "Frequency in description matches metadata (weekly, monthly)"

# --- Keyword Matches sheet ---
$ws = $wb.Worksheets.Item("Keyword Matches")

# Column C width: 48.4 -> 29.7 (closest achievable via pixel-quantized ColumnWidth is 29.666667)
$ws.Columns.Item(3).ColumnWidth = 28.8

$ws.Range("C5").Value = "by the infosec team"
$ws.Range("C6").Value = "by the finance team bef"
$ws.Range("C12").Value = "on a monthly basis"
$ws.Range("C14").Value = "quarterly"
$ws.Range("C15").Value = "as needed"
$ws.Range("C20").Value = "None"
$ws.Range("C23").Value = "None"
$ws.Range("C27").Value = "daily"
$ws.Range("C29").Value = "ad hoc, on a weekly basis"

# --- Enhancement Feedback sheet ---
$ws = $wb.Worksheets.Item("Enhancement Feedback")

$ws.Range("C4").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Monthly)"
$ws.Range("C5").Value = "Align the frequency in the description with the declared frequency (Ongoing)"
$ws.Range("C6").Value = "Align the frequency in the description with the declared frequency (Daily)"
$ws.Range("C7").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Monthly)"
$ws.Range("C8").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (As needed)"
$ws.Range("C9").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Quarterly)"
$ws.Range("C10").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Weekly)"
$ws.Range("C11").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Monthly)"
$ws.Range("C15").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Replace vague timing term 'as needed' with a specific timeframe or frequency.; Align the frequency in the description with the declared frequency (As needed)"
$ws.Range("C16").Value = "Align the frequency in the description with the declared frequency (Within 3 days)"
$ws.Range("C17").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Ongoing)"
$ws.Range("C18").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Regulatory)"
$ws.Range("C19").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Ongoing)"
$ws.Range("C20").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Prompt)"
$ws.Range("C21").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Internal)"
$ws.Range("C23").Value = "No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days)."
$ws.Range("C29").Value = "While 'ad-hoc' is an allowed frequency, the control would be stronger if it specified what triggers the ad-hoc review.; Multiple frequencies detected. Consider whether this is describing a process rather than a single control."

# --- Executive Summary sheet ---
$ws = $wb.Worksheets.Item("Executive Summary")

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "42.3"
$ws.Range("B4").NumberFormat = "General"
$ws.Range("B15").Value = "17 (56.7%)"
$ws.Range("B24").Value = "15 (50.0%)"
$ws.Range("B25").Value = "15 (50.0%)"
